# Clean database pass: convert the "roboticS1Prep" column (I) from the
# free-text "No" values into a proper boolean FALSE with a TRUE/FALSE
# custom number format, for every data row (2-27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 27

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("I$r")
    $cell.Value = $false
    $cell.NumberFormat = '"TRUE";"TRUE";"FALSE"'
}

# Selection follows the edited column, matching the author's last click.
$ws.Range("I2:I27").Select() | Out-Null
